$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the week's progress row (row 17, date 16/5/2025 already present in D17)
$ws.Range("E17").Value = 276
$ws.Range("F17").Value = 80
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 650
$ws.Range("J17").Value = "N/A"

# Update the last selected cell as recorded in the saved file
$ws.Range("G26").Select()
